$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet
$ws.Name = "TestingPackage"

# Fill in the test-case data grid (rows 2-11)
$ws.Range("A2").Value = "TestScenario_1"
$ws.Range("B2").Value = "TestScenario_1.TestCase_1"
$ws.Range("C2").Value = "New Account"
$ws.Range("D2").Value = "User Needs to Login to Salesforce, from the browser with correct credentials"
$ws.Range("F2").Value = "Step 1"
$ws.Range("G2").Value = "Click on the Account tab, and click on New button"
$ws.Range("H2").Value = "User should be navigated to the New  Account Page"
$ws.Range("E2").ClearContents()
$ws.Range("I2").ClearContents()
$ws.Range("J2").ClearContents()
$ws.Range("E3").Value = "Valid value for required field Account Name "
$ws.Range("F3").Value = "Step 2"
$ws.Range("G3").Value = "Input valid value in the  Account Name field."
$ws.Range("H3").Value = "User should be able to input value for the Account Name field."
$ws.Range("E4").Value = "Valid value for required field  "
$ws.Range("F4").Value = "Step 3"
$ws.Range("G4").Value = "Input valid value in the   field."
$ws.Range("H4").Value = "User should be able to input value for the  field."
$ws.Range("E5").Value = "Valid value for required field Annual Revenue, value should be greater than 50000 to submit for Approval towards to the assigned approver "
$ws.Range("F5").Value = "Step 4"
$ws.Range("G5").Value = "Input valid value in the  Annual Revenue field."
$ws.Range("H5").Value = "Value accepted for Annual Revenue field."
$ws.Range("F6").Value = "Step 5"
$ws.Range("G6").Value = "Click on Save button to save Account with fields"
$ws.Range("H6").Value = "User should be able to validate that a New Account is created"
$ws.Range("F7").Value = "Step 6"
$ws.Range("G7").Value = "On the Account page, Click on 'Submit for Approval' button to Submit for Approval."
$ws.Range("H7").Value = "Pop-up confirming to submit the record for Approval is displayed."
$ws.Range("F8").Value = "Step 7"
$ws.Range("G8").Value = "Click on 'OK' button to submit the record for Approval."
$ws.Range("H8").Value = "The record will be displayed under Approval History section with the status 'Pending'."
$ws.Range("F9").Value = "Step 8"
$ws.Range("G9").Value = "If the user is navigated to the 'Choose Approver' page, he should be able to input the Next Approver."
$ws.Range("H9").Value = "User is able to input the Next Approver."
$ws.Range("F10").Value = "Step 9"
$ws.Range("G10").Value = "Click on 'Send to Next Approver' button."
$ws.Range("H10").Value = "The record will be displayed under Approval History section with the status 'Pending' and the updated Approver Name."
$ws.Range("F11").Value = "Step 10"
$ws.Range("G11").Value = "For this Approval process, if e-mail notification is configured, validate that the Approver receives an e-mail with Approve/Reject request."
$ws.Range("H11").Value = "User should be able to view the e-mail to either Approve/Reject the request."

# Resize the Table1 list object + column widths
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:J11"))

$ws.Columns.Item(2).ColumnWidth = 25.5
$ws.Columns.Item(3).ColumnWidth = 13.833333333333332
$ws.Columns.Item(4).ColumnWidth = 69.16666666666667
$ws.Columns.Item(5).ColumnWidth = 124.66666666666667
$ws.Columns.Item(6).ColumnWidth = 8.5
$ws.Columns.Item(7).ColumnWidth = 122.33333333333334
$ws.Columns.Item(8).ColumnWidth = 106.66666666666667
